$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the weekly price data between row 2 and row 3 (keep D as raw date serials)
$ws.Range("D2").Value = 44749
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 17000
$ws.Range("L2").Value = 18000
$ws.Range("M2").Value = 17556
$ws.Range("P2").Value = 1170

$ws.Range("D3").Value = 44839
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 16000
$ws.Range("M3").Value = 15600
$ws.Range("P3").Value = 1040
